$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# Update the label in C2 ("板块轮动策略" -> "板块热点轮动策略")
$ws.Cells.Item(2, 3).Value = "板块热点轮动策略"

# Remove the stray leftover rows (3-7, 11, 14, 16) below the real data,
# leaving only the header row and the single data row.
$ws.Range("A3:J16").Clear()

# Widen columns B and C to fit the (now longer) strategy-name text.
$ws.Columns.Item(2).ColumnWidth = 15.714285714285714
$ws.Columns.Item(3).ColumnWidth = 14.571428571428571

# Move the active selection to C6.
$ws.Range("C6").Select()
